# StyleTag InputData.xlsx - "Login" sheet rework:
#  - add a new "two blankspaces in each cell" negative-scenario row
#  - split the previously rich-text "test3456@styletag" + "." cell into one plain run
#  - add reason/explanation notes in column C for each negative scenario
#  - add a second "LoginButton enabled" (bold, blue) section with valid-login scenarios
#  - re-point / re-create the mailto hyperlinks to match the re-shuffled rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- wipe the old hyperlinks + cell contents so we can rebuild cleanly ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Cells.Clear()

# --- row 1-2: headers ---
$ws.Range("A1").Value = "LoginButton disability"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Value = "Login ID"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = "Password"
$ws.Range("B2").Font.Bold = $true

# --- row 3: both fields blank (two spaces each) ---
$ws.Range("A3").Value = "  "
$ws.Range("B3").Value = "  "
$ws.Range("C3").Value = "two blankspaces in each cell (both the fields are empty)"
$ws.Range("C3").Font.Bold = $true

# --- row 4 ---
$ws.Range("A4").Value = "test3456"
$ws.Range("B4").Value = "styletag123"
$ws.Range("B4").Font.Color = 16711680

# --- row 5 ---
$ws.Range("A5").Value = "test3456@"
$ws.Range("B5").Value = "styletag123"
$ws.Range("B5").Font.Color = 16711680

# --- row 6 (mailto hyperlink) ---
$ws.Range("A6").Value = "test3456@st"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:test3456@st", $null, $null, "test3456@st")
$ws.Range("B6").Value = "styletag123"
$ws.Range("B6").Font.Color = 16711680

# --- row 7 (mailto hyperlink) ---
$ws.Range("A7").Value = "test3456@styletag."
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:test3456@styletag.", $null, $null, "test3456@styletag.")
$ws.Range("B7").Value = "styletag123"
$ws.Range("B7").Font.Color = 16711680

# --- row 8 (mailto hyperlink) ---
$ws.Range("A8").Value = "test3456@styletag.c"
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:test3456@styletag.c", $null, $null, "test3456@styletag.c")
$ws.Range("B8").Value = "styletag123"
$ws.Range("B8").Font.Color = 16711680

# --- row 9: blank email ---
$ws.Range("A9").Value = "      "
$ws.Range("B9").Value = "styletag123"
$ws.Range("C9").Value = "blank emailis"
$ws.Range("C9").Font.Bold = $true

# --- row 10: blank password (mailto hyperlink) ---
$ws.Range("A10").Value = "test3456@styletag.com"
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:test3456@styletag.com", $null, $null, "test3456@styletag.com")
$ws.Range("B10").Value = "        "
$ws.Range("C10").Value = "blank passwors"
$ws.Range("C10").Font.Bold = $true

# --- row 11 ---
$ws.Range("A11").Value = "Test3456@.com"
$ws.Range("B11").Value = "styletag123"
$ws.Range("B11").Font.Color = 16711680

# --- row 12 ---
$ws.Range("A12").Value = "@.com"
$ws.Range("B12").Value = "styletag123"
$ws.Range("B12").Font.Color = 16711680

# --- row 13: second header block (bold + blue) ---
$ws.Range("A13").Value = "LoginButton enabled"
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Font.Color = 16711680
$ws.Range("B13").Value = "following data are not valid for login"
$ws.Range("B13").Font.Bold = $true

# --- row 14: invalid password (mailto hyperlink) ---
$ws.Range("A14").Value = "test3456@styletag.com"
$ws.Hyperlinks.Add($ws.Range("A14"), "mailto:test3456@styletag.com", $null, $null, "test3456@styletag.com")
$ws.Range("B14").Value = "styletag12"
$ws.Range("C14").Value = "invalid password"
$ws.Range("C14").Font.Bold = $true

# --- row 15: not registered credentials (mailto hyperlink) ---
$ws.Range("A15").Value = "test34567@styletag.com"
$ws.Hyperlinks.Add($ws.Range("A15"), "mailto:test34567@styletag.com", $null, $null, "test34567@styletag.com")
$ws.Range("B15").Value = "styletag123"
$ws.Range("C15").Value = "not registered credentials"
$ws.Range("C15").Font.Bold = $true

# --- row 16: valid data header ---
$ws.Range("A16").Value = "Valid data"
$ws.Range("A16").Font.Bold = $true

# --- row 17: registered, valid data (mailto hyperlink) ---
$ws.Range("A17").Value = "test3456@styletag.com"
$ws.Hyperlinks.Add($ws.Range("A17"), "mailto:test3456@styletag.com", $null, $null, "test3456@styletag.com")
$ws.Range("B17").Value = "styletag123"
$ws.Range("C17").Value = "Registered data"
$ws.Range("C17").Font.Bold = $true

$ws.Range("C9").Select()
